# Insert a new data row at row 391 (pushing the existing rows 391-414 down
# to 392-415) and populate it with the new "Vega Modelo de Temuco - Coliflor"
# price record. This also grows the sheet's used range from R414 to R415.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 391..414 down to 392..415, carrying formatting (incl. the date
# number format on column D) from the row being pushed down - same as
# Excel's native "Insert Sheet Rows" command.
$ws.Rows.Item(391).Insert()

# Fill in the newly inserted row with the new record's data.
$ws.Range("A391").Value = 10
$ws.Range("B391").Value = "Vega Modelo de Temuco"
$ws.Range("C391").Value = "La Araucanía"
$ws.Range("D391").Value = 44753
$ws.Range("E391").Value = 9
$ws.Range("F391").Value = 100112008
$ws.Range("G391").Value = "Coliflor"
$ws.Range("H391").Value = "Sin especificar"
$ws.Range("I391").Value = "Primera"
$ws.Range("J391").Value = 300
$ws.Range("K391").Value = 1200
$ws.Range("L391").Value = 1200
$ws.Range("M391").Value = 1200
$ws.Range("N391").Value = "$/unidad"
$ws.Range("O391").Value = "Región Metropolitana"
$ws.Range("P391").Value = 1200
$ws.Range("Q391").Value = 1
$ws.Range("R391").Value = "Hortaliza"
